$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Companies")

# Row 3 updates: Opportunity/Engagement refresh
$ws.Range("C3").Value = "Project Miura"
$ws.Range("D3").Value = "120228"
$ws.Range("F3").Value = "Project Chavez"
$ws.Range("G3").Value = "120227"
$ws.Range("H3").Value = "Book Prep/Underwriting"

$ws.Range("E8").Select()
